# Actualización automática 2025-10-08 14:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M14").Value = 1071.28
$ws1.Range("P14").Value = 143.08
$ws1.Range("M48").Value = 111.13
$ws1.Range("M59").Value = "3 de 57"
$ws1.Range("P59").Value = "2 de 57"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F14").Value = 2122.39
$ws2.Range("F48").Value = 501.5
$ws2.Range("F59").Value = 4505.73

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D8").Value = 226.58
$ws3.Range("E8").Value = 222.22162917203
$ws3.Range("F8").Value = 0.5048555648472249

$ws3.Range("D12").Value = 2221.89
$ws3.Range("E12").Value = 46402.17
$ws3.Range("F12").Value = 0.04569527925064258

$ws3.Range("D14").Value = 4562.05
$ws3.Range("E14").Value = 95335.94284188785
$ws3.Range("F14").Value = 0.04566708369426923
